# "Changes of New Pre-Prod URL" - refresh the FedEx shipment tracking
# numbers in column P ("ShipmentTracking") for every data row (2-26) to
# the new tracking numbers generated against the new pre-prod endpoint.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTracking = @(
    "320018589548",
    "320018589559",
    "320018589581",
    "320018589607",
    "320018589640",
    "320018589662",
    "320018589695",
    "320018589710",
    "320018589743",
    "320018589765",
    "320018589802",
    "320018589824",
    "320018589857",
    "320018589879",
    "320018589905",
    "320018589927",
    "320018589960",
    "320018589982",
    "320018590015",
    "320018590037",
    "320018590060",
    "320018590070",
    "320018590081",
    "320018590092",
    "320018590107"
)

for ($i = 0; $i -lt $newTracking.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("P$row")

    # Tracking numbers are all-digit strings; format the cell as Text
    # first so Excel keeps the new value as a string (matching the
    # original t="s" shared-string cells) instead of coercing it to a
    # number, then restore the default "Normal" style so no stray
    # number-format is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $newTracking[$i]
    $cell.Style = "Normal"
}
